$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 339.3
$ws.Range("I2").Value = 48.875
$ws.Range("K2").Value = 48.875
$ws.Range("M2").Value = 64.125
$ws.Range("H18").Value = 590
$ws.Range("I18").Value = 590
$ws.Range("K18").Value = 590
$ws.Range("M18").Value = -306
$ws.Range("H43").Value = 3837.7693
$ws.Range("J43").Value = 4199.1
$ws.Range("L43").Value = 4199.1
$ws.Range("N43").Value = -4337.1
$ws.Range("H62").Value = 5722379.5
$ws.Range("I62").Value = 6501204
$ws.Range("K62").Value = 6501204
$ws.Range("M62").Value = -6500580
$ws.Range("H65").Value = 5722379.5
$ws.Range("I65").Value = 6501204
$ws.Range("K65").Value = 32506020
$ws.Range("M65").Value = -32502900
$ws.Range("H76").Value = 166672020
$ws.Range("I76").Value = 250004270
$ws.Range("J76").Value = 7502
$ws.Range("K76").Value = 250004270
$ws.Range("L76").Value = 7502
$ws.Range("M76").Value = -250003955
$ws.Range("N76").Value = -8132
$ws.Range("H79").Value = 166672020
$ws.Range("I79").Value = 250004270
$ws.Range("J79").Value = 7502
$ws.Range("K79").Value = 250004270
$ws.Range("L79").Value = 7502
$ws.Range("M79").Value = -250003178
$ws.Range("N79").Value = -9686
$ws.Range("H98").Value = 772
$ws.Range("I98").Value = 803.0833
$ws.Range("J98").Value = 26
$ws.Range("K98").Value = 803.0833
$ws.Range("L98").Value = 26
$ws.Range("M98").Value = 694.9167
$ws.Range("N98").Value = -3022
$ws.Range("H106").Value = 2661
$ws.Range("I106").Value = 3337.4614
$ws.Range("J106").Value = 1683.8889
$ws.Range("K106").Value = 3337.4614
$ws.Range("L106").Value = 1683.8889
$ws.Range("M106").Value = -2706.4614
$ws.Range("N106").Value = -2945.8889
$ws.Range("H107").Value = 2867.1052
$ws.Range("I107").Value = 1104.3636
$ws.Range("J107").Value = 5290.875
$ws.Range("K107").Value = 1104.3636
$ws.Range("L107").Value = 5290.875
$ws.Range("M107").Value = 815.6364000000001
$ws.Range("N107").Value = -9130.875
$ws.Range("H122").Value = 772
$ws.Range("I122").Value = 803.0833
$ws.Range("J122").Value = 26
$ws.Range("K122").Value = 2409.2499
$ws.Range("L122").Value = 78
$ws.Range("M122").Value = 40.7501000000002
$ws.Range("N122").Value = -4978
$ws.Range("H132").Value = 287983.97
$ws.Range("I132").Value = 477510.22
$ws.Range("K132").Value = 1432530.66
$ws.Range("M132").Value = -1430000.66
$ws.Range("H137").Value = 4310.5
$ws.Range("J137").Value = 7011.5
$ws.Range("L137").Value = 21034.5
$ws.Range("N137").Value = -26134.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 473.3
$ws.Range("I4").Value = 450.625
$ws.Range("K4").Value = 450.625
$ws.Range("M4").Value = -334.625
$ws.Range("I32").Value = 671.2
$ws.Range("K32").Value = 671.2
$ws.Range("M32").Value = -384.2
$ws.Range("H34").Value = 18512
$ws.Range("I34").Value = 18512
$ws.Range("K34").Value = 18512
$ws.Range("M34").Value = -18241
$ws.Range("H40").Value = 3999.5
$ws.Range("I40").Value = 3999.5
$ws.Range("K40").Value = 3999.5
$ws.Range("M40").Value = -3823.5
$ws.Range("H45").Value = 1550.5264
$ws.Range("I45").Value = 1594.9333
$ws.Range("J45").Value = 1384
$ws.Range("K45").Value = 1594.9333
$ws.Range("L45").Value = 1384
$ws.Range("M45").Value = -1217.9333
$ws.Range("N45").Value = -2138
$ws.Range("H61").Value = 3834.5715
$ws.Range("I61").Value = 1986.8
$ws.Range("K61").Value = 1986.8
$ws.Range("M61").Value = -1774.8
$ws.Range("H76").Value = 60999.5
$ws.Range("J76").Value = 60999.5
$ws.Range("L76").Value = 60999.5
$ws.Range("N76").Value = -61675.5
$ws.Range("H79").Value = 60999.5
$ws.Range("J79").Value = 60999.5
$ws.Range("L79").Value = 60999.5
$ws.Range("N79").Value = -63339.5
$ws.Range("H88").Value = 5863.143
$ws.Range("I88").Value = 1928.1428
$ws.Range("K88").Value = 1928.1428
$ws.Range("M88").Value = -1522.1428
$ws.Range("H91").Value = 5863.143
$ws.Range("I91").Value = 1928.1428
$ws.Range("K91").Value = 1928.1428
$ws.Range("M91").Value = -524.1428000000001
$ws.Range("H97").Value = 738.087
$ws.Range("I97").Value = 568.8946999999999
$ws.Range("K97").Value = 568.8946999999999
$ws.Range("M97").Value = -72.89469999999994
$ws.Range("H110").Value = 1535
$ws.Range("I110").Value = 1537.5
$ws.Range("J110").Value = 1500
$ws.Range("K110").Value = 1537.5
$ws.Range("L110").Value = 1500
$ws.Range("M110").Value = 507.5
$ws.Range("N110").Value = -5590
$ws.Range("H122").Value = 2556
$ws.Range("I122").Value = 1417.8182
$ws.Range("K122").Value = 4253.4546
$ws.Range("M122").Value = -1803.4546
$ws.Range("H132").Value = 1252283.6
$ws.Range("I132").Value = 1686283
$ws.Range("K132").Value = 5058849
$ws.Range("M132").Value = -5056319
$ws.Range("H136").Value = 3834.5715
$ws.Range("I136").Value = 1986.8
$ws.Range("K136").Value = 5960.4
$ws.Range("M136").Value = -3410.4
$ws.Range("H138").Value = 90000
$ws.Range("J138").Value = 90000
$ws.Range("L138").Value = 90000
$ws.Range("N138").Value = -100280
$ws.Range("H139").Value = 900715
$ws.Range("J139").Value = 900715
$ws.Range("L139").Value = 900715
$ws.Range("N139").Value = -910995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10206831
$ws.Range("I20").Value = 23812804
$ws.Range("K20").Value = 23812804
$ws.Range("M20").Value = -23812557
$ws.Range("H94").Value = 2509.842
$ws.Range("I94").Value = 1176.7693
$ws.Range("J94").Value = 5398.1665
$ws.Range("K94").Value = 1176.7693
$ws.Range("L94").Value = 5398.1665
$ws.Range("M94").Value = -725.7692999999999
$ws.Range("N94").Value = -6300.1665
$ws.Range("H99").Value = 7791.081
$ws.Range("I99").Value = 7791.081
$ws.Range("K99").Value = 7791.081
$ws.Range("M99").Value = -6293.081
$ws.Range("H105").Value = 4910.8887
$ws.Range("I105").Value = 5343.5
$ws.Range("J105").Value = 1450
$ws.Range("K105").Value = 5343.5
$ws.Range("L105").Value = 1450
$ws.Range("M105").Value = -3596.5
$ws.Range("N105").Value = -4944
$ws.Range("H107").Value = 9092615
$ws.Range("I107").Value = 12501343
$ws.Range("J107").Value = 2674.3333
$ws.Range("K107").Value = 12501343
$ws.Range("L107").Value = 2674.3333
$ws.Range("M107").Value = -12499423
$ws.Range("N107").Value = -6514.3333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 55560650
$ws.Range("I16").Value = 100003460
$ws.Range("J16").Value = 7138.75
$ws.Range("K16").Value = 100003460
$ws.Range("L16").Value = 7138.75
$ws.Range("M16").Value = -100003173
$ws.Range("N16").Value = -7712.75
$ws.Range("H31").Value = 2132.074
$ws.Range("I31").Value = 1035.5834
$ws.Range("J31").Value = 3009.2666
$ws.Range("K31").Value = 1035.5834
$ws.Range("L31").Value = 3009.2666
$ws.Range("M31").Value = -740.5834
$ws.Range("N31").Value = -3599.2666
$ws.Range("H34").Value = 2132.074
$ws.Range("I34").Value = 1035.5834
$ws.Range("J34").Value = 3009.2666
$ws.Range("K34").Value = 1035.5834
$ws.Range("L34").Value = 3009.2666
$ws.Range("M34").Value = -833.5834
$ws.Range("N34").Value = -3413.2666
$ws.Range("H105").Value = 111115650
$ws.Range("I105").Value = 166669310
$ws.Range("J105").Value = 8333.333000000001
$ws.Range("K105").Value = 166669310
$ws.Range("L105").Value = 8333.333000000001
$ws.Range("M105").Value = -166667563
$ws.Range("N105").Value = -11827.333
$ws.Range("H107").Value = 1415
$ws.Range("I107").Value = 1274.2858
$ws.Range("J107").Value = 1579.1666
$ws.Range("K107").Value = 1274.2858
$ws.Range("L107").Value = 1579.1666
$ws.Range("M107").Value = 645.7141999999999
$ws.Range("N107").Value = -5419.1666
$ws.Range("H113").Value = 55560650
$ws.Range("I113").Value = 100003460
$ws.Range("J113").Value = 7138.75
$ws.Range("K113").Value = 100003460
$ws.Range("L113").Value = 7138.75
$ws.Range("M113").Value = -100001290
$ws.Range("N113").Value = -11478.75
$ws.Range("H122").Value = 1965.8572
$ws.Range("J122").Value = 1944
$ws.Range("L122").Value = 5832
$ws.Range("N122").Value = -10732
$ws.Range("H131").Value = 44000.5
$ws.Range("J131").Value = 44000.5
$ws.Range("L131").Value = 44000.5
$ws.Range("N131").Value = -54080.5
$ws.Range("H132").Value = 3537.725
$ws.Range("I132").Value = 2647.6128
$ws.Range("J132").Value = 6603.6665
$ws.Range("K132").Value = 7942.8384
$ws.Range("L132").Value = 19810.9995
$ws.Range("M132").Value = -5412.8384
$ws.Range("N132").Value = -24870.9995
$ws.Range("H134").Value = 27032468
$ws.Range("I134").Value = 76924890
$ws.Range("K134").Value = 230774670
$ws.Range("M134").Value = -230772135

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 10246990
$ws.Range("J4").Value = 1234567
$ws.Range("L4").Value = 3703701
$ws.Range("N4").Value = -3703925
$ws.Range("H7").Value = 613.1111
$ws.Range("I7").Value = 456.33334
$ws.Range("J7").Value = 926.6667
$ws.Range("K7").Value = 1369.00002
$ws.Range("L7").Value = 2780.0001
$ws.Range("M7").Value = -1257.00002
$ws.Range("N7").Value = -3004.0001
$ws.Range("H14").Value = 1576.875
$ws.Range("I14").Value = 1576.875
$ws.Range("K14").Value = 4730.625
$ws.Range("M14").Value = -4557.625
$ws.Range("H15").Value = 64.5
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H87").Value = 20702.857
$ws.Range("I87").Value = 5306.6665
$ws.Range("K87").Value = 15919.9995
$ws.Range("M87").Value = -14671.9995
$ws.Range("H90").Value = 20702.857
$ws.Range("I90").Value = 5306.6665
$ws.Range("K90").Value = 47759.9985
$ws.Range("M90").Value = -41519.9985
$ws.Range("H111").Value = 6204.5
$ws.Range("J111").Value = 11624.75
$ws.Range("L111").Value = 34874.25
$ws.Range("N111").Value = -41008.25
$ws.Range("H122").Value = 193152.92
$ws.Range("J122").Value = 209210.3
$ws.Range("L122").Value = 1882892.7
$ws.Range("N122").Value = -1887792.7
$ws.Range("H140").Value = 75762320
$ws.Range("I140").Value = 104171310
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 312513930
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = -312508750
$ws.Range("N140").Value = -25360

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 66666
$ws.Range("I21").Value = 66666
$ws.Range("K21").Value = 66666
$ws.Range("M21").Value = -66493
$ws.Range("H24").Value = 2005400
$ws.Range("I24").Value = 2005400
$ws.Range("K24").Value = 2005400
$ws.Range("M24").Value = -2005227
$ws.Range("H30").Value = 66666
$ws.Range("I30").Value = 66666
$ws.Range("K30").Value = 66666
$ws.Range("M30").Value = -66561
$ws.Range("H70").Value = 10422.223
$ws.Range("I70").Value = 10147.25
$ws.Range("J70").Value = 10642.2
$ws.Range("K70").Value = 10147.25
$ws.Range("L70").Value = 10642.2
$ws.Range("M70").Value = -9877.25
$ws.Range("N70").Value = -11182.2
$ws.Range("H73").Value = 10422.223
$ws.Range("I73").Value = 10147.25
$ws.Range("J73").Value = 10642.2
$ws.Range("K73").Value = 10147.25
$ws.Range("L73").Value = 10642.2
$ws.Range("M73").Value = -9211.25
$ws.Range("N73").Value = -12514.2
$ws.Range("H80").Value = 3642.2307
$ws.Range("I80").Value = 1598.75
$ws.Range("J80").Value = 6911.8
$ws.Range("K80").Value = 1598.75
$ws.Range("L80").Value = 6911.8
$ws.Range("M80").Value = -600.75
$ws.Range("N80").Value = -8907.799999999999
$ws.Range("H83").Value = 3642.2307
$ws.Range("I83").Value = 1598.75
$ws.Range("J83").Value = 6911.8
$ws.Range("K83").Value = 7993.75
$ws.Range("L83").Value = 34559
$ws.Range("M83").Value = -3001.75
$ws.Range("N83").Value = -44543
$ws.Range("H122").Value = 8666.639999999999
$ws.Range("I122").Value = 7897.75
$ws.Range("K122").Value = 23693.25
$ws.Range("M122").Value = -21243.25
$ws.Range("H126").Value = 20007786
$ws.Range("J126").Value = 19162.5
$ws.Range("L126").Value = 57487.5
$ws.Range("N126").Value = -62427.5
$ws.Range("H132").Value = 41670144
$ws.Range("I132").Value = 66669460
$ws.Range("J132").Value = 4610.5557
$ws.Range("K132").Value = 200008380
$ws.Range("L132").Value = 13831.6671
$ws.Range("M132").Value = -200005850
$ws.Range("N132").Value = -18891.6671

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 16981.818
$ws.Range("I20").Value = 10850
$ws.Range("K20").Value = 10850
$ws.Range("M20").Value = -10624
$ws.Range("H22").Value = 1119.6
$ws.Range("I22").Value = 1166.6666
$ws.Range("K22").Value = 1166.6666
$ws.Range("M22").Value = -871.6666
$ws.Range("H27").Value = 1119.6
$ws.Range("I27").Value = 1166.6666
$ws.Range("K27").Value = 1166.6666
$ws.Range("M27").Value = -1059.6666
$ws.Range("H46").Value = 23811380
$ws.Range("J46").Value = 50003000
$ws.Range("L46").Value = 50003000
$ws.Range("N46").Value = -50003376
$ws.Range("H68").Value = 1566.6666
$ws.Range("I68").Value = 1566.6666
$ws.Range("K68").Value = 1566.6666
$ws.Range("M68").Value = -817.6666
$ws.Range("H71").Value = 1566.6666
$ws.Range("I71").Value = 1566.6666
$ws.Range("K71").Value = 7833.333000000001
$ws.Range("M71").Value = -4089.333000000001
$ws.Range("H88").Value = 24390.75
$ws.Range("I88").Value = 24188
$ws.Range("J88").Value = 24593.5
$ws.Range("K88").Value = 24188
$ws.Range("L88").Value = 24593.5
$ws.Range("M88").Value = -23760
$ws.Range("N88").Value = -25449.5
$ws.Range("H91").Value = 24390.75
$ws.Range("I91").Value = 24188
$ws.Range("J91").Value = 24593.5
$ws.Range("K91").Value = 24188
$ws.Range("L91").Value = 24593.5
$ws.Range("M91").Value = -22706
$ws.Range("N91").Value = -27557.5
$ws.Range("H124").Value = 39999.5
$ws.Range("J124").Value = 39999.5
$ws.Range("L124").Value = 39999.5
$ws.Range("N124").Value = -49819.5
$ws.Range("H131").Value = 56216
$ws.Range("J131").Value = 54246.855
$ws.Range("L131").Value = 54246.855
$ws.Range("N131").Value = -64326.855
$ws.Range("H132").Value = 3030
$ws.Range("J132").Value = 4162.25
$ws.Range("L132").Value = 12486.75
$ws.Range("N132").Value = -17546.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 45000
$ws.Range("J70").Value = 45000
$ws.Range("L70").Value = 45000
$ws.Range("N70").Value = -45630
$ws.Range("H73").Value = 45000
$ws.Range("J73").Value = 45000
$ws.Range("L73").Value = 45000
$ws.Range("N73").Value = -47184
$ws.Range("H81").Value = 765.6
$ws.Range("I81").Value = 765.6
$ws.Range("K81").Value = 1531.2
$ws.Range("M81").Value = -470.2
$ws.Range("H84").Value = 765.6
$ws.Range("I84").Value = 765.6
$ws.Range("K84").Value = 7656
$ws.Range("M84").Value = -2352
$ws.Range("H113").Value = 44218320
$ws.Range("I113").Value = 51587876
$ws.Range("K113").Value = 154763628
$ws.Range("M113").Value = -154761458
